$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 odds updates per diff
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2.7
$ws.Range("I10").Value = 3.4

$ws.Range("W10").Value = 6
$ws.Range("X10").Value = 10
$ws.Range("Y10").Value = 11
$ws.Range("Z10").Value = 23
$ws.Range("AA10").Value = 26

$ws.Range("AE10").Value = 19
$ws.Range("AF10").Value = 81
$ws.Range("AH10").Value = 15

$ws.Range("BA10").Value = 126
